$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Fix the horizontal-line VML shape (w:pict / v:shape) so it uses
#    the "modern" Word export form: adds alt="", mso-wrap-edited/
#    mso-width-percent/mso-height-percent to style, switches from a
#    coordorigin-based path to a simple 0-based path with proper
#    coordsize + connector locations.
# ------------------------------------------------------------------

$shapeParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $cr = $cand.Range
    if ($cr.Start -eq $cr.End) {
        # zero-length (no visible text) paragraph; the drawing-only
        # paragraph that holds our v:shape lives around here. Narrow
        # it down by checking neighbouring text (the shape paragraph
        # sits immediately before the "${id_usuario}" paragraph).
        $nextPara = $d.Paragraphs.Item($i + 1)
        if ($nextPara.Range.Text -like "*id_usuario*") {
            $shapeParaIndex = $i
        }
    }
}

if ($shapeParaIndex -ne -1) {
    $shapeRange = $d.Paragraphs.Item($shapeParaIndex).Range
    $shapeXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="0FD1D374" w14:textId="77777777" w:rsidR="009B33D1" w:rsidRPr="009F7B04" w:rsidRDefault="00000000"><w:pPr><w:pStyle w:val="Textoindependiente"/><w:spacing w:before="4"/><w:rPr><w:sz w:val="14"/></w:rPr></w:pPr><w:r><w:pict w14:anchorId="0FD1D388"><v:shape id="_x0000_s1026" alt="" style="position:absolute;margin-left:220.4pt;margin-top:11.1pt;width:149.25pt;height:.1pt;z-index:-251658752;mso-wrap-edited:f;mso-width-percent:0;mso-height-percent:0;mso-wrap-distance-left:0;mso-wrap-distance-right:0;mso-position-horizontal-relative:page;mso-width-percent:0;mso-height-percent:0" coordsize="2985,1270" path="m,l2985,e" filled="f"><v:path arrowok="t" o:connecttype="custom" o:connectlocs="0,0;1895475,0" o:connectangles="0,0"/><w10:wrap type="topAndBottom" anchorx="page"/></v:shape></w:pict></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
    $shapeRange.InsertXML($shapeXml)
}

# ------------------------------------------------------------------
# 2) Collapse the "30 DE JUNIO DE 2023" date run-cluster: drop the
#    yellow highlight (w:shd) from every run and merge the separately
#    spaced runs back into fewer, simpler runs.
# ------------------------------------------------------------------

$dateParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -eq " 30 DE JUNIO DE 2023 ") {
        $dateParaIndex = $i
    }
}

if ($dateParaIndex -ne -1) {
    $dateRange = $d.Paragraphs.Item($dateParaIndex).Range
    $dateXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="0FD1D384" w14:textId="40C73262" w:rsidR="009B33D1" w:rsidRDefault="0098720A"><w:pPr><w:pStyle w:val="TableParagraph"/><w:spacing w:before="1" w:line="223" w:lineRule="exact"/><w:ind w:left="193"/><w:rPr><w:b/><w:sz w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="009F7B04"><w:rPr><w:b/><w:sz w:val="20"/></w:rPr><w:t>30</w:t></w:r><w:r w:rsidR="00A377B5" w:rsidRPr="009F7B04"><w:rPr><w:b/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve"> DE </w:t></w:r><w:r w:rsidRPr="009F7B04"><w:rPr><w:b/><w:sz w:val="20"/></w:rPr><w:t>JUNIO</w:t></w:r><w:r w:rsidR="00A377B5" w:rsidRPr="009F7B04"><w:rPr><w:b/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve"> DE 202</w:t></w:r><w:r w:rsidRPr="009F7B04"><w:rPr><w:b/><w:sz w:val="20"/></w:rPr><w:t>3</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
    $dateRange.InsertXML($dateXml)
}
